# Refresh the cryptocurrency price / 1h-volume snapshot (Sat Sep 14 2024 run).
#
# Every data cell on the sheet is stored as literal text (t="inlineStr" in the
# OOXML), including the "Price" column, which often looks numeric (e.g. "550.96").
# Assigning such a string straight to Range.Value would make Excel silently
# reinterpret it as a number, so those assignments are prefixed with a leading
# apostrophe (the normal Excel text qualifier) to keep them as text, exactly as
# the source data has them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.832.40'
$ws.Range("E2").Value = '  +2.58%  '

$ws.Range("D3").Value = '2.413.30'
$ws.Range("E3").Value = '  +1.95%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").Value = '''550.96'

$ws.Range("D6").Value = '''137.21'
$ws.Range("E6").Value = '  +2.35%  '

$ws.Range("D7").Value = '''0.999'
$ws.Range("E7").Value = '  -0.03%  '

$ws.Range("E8").Value = '  +2.67%  '

$ws.Range("E9").Value = '  -0.58%  '

$ws.Range("D10").Value = '''5.73'
$ws.Range("E10").Value = '  +3.58%  '

$ws.Range("E11").Value = '  -1.94%  '

$ws.Range("D12").Value = '''0.355'
$ws.Range("E12").Value = '  -0.34%  '

$ws.Range("D13").Value = '''24.72'
$ws.Range("E13").Value = '  +2.74%  '

$ws.Range("D14").Value = '2.842.03'
$ws.Range("E14").Value = '  +2.07%  '

$ws.Range("D15").Value = '59.801.18'
$ws.Range("E15").Value = '  +2.63%  '

$ws.Range("D17").Value = '2.408.47'
$ws.Range("E17").Value = '  +2.48%  '

$ws.Range("D18").Value = '''11.28'
$ws.Range("E18").Value = '  +2.87%  '

$ws.Range("E19").Value = '  +0.85%  '

$ws.Range("D20").Value = '''330.61'
$ws.Range("E20").Value = '  -0.23%  '

$ws.Range("D21").Value = '''6.69'
$ws.Range("E21").Value = '  -2.55%  '

$ws.Range("E22").Value = '  -0.10%  '

$ws.Range("D23").Value = '''65.79'
$ws.Range("E23").Value = '  +3.45%  '

$ws.Range("E24").Value = '  +2.73%  '

$ws.Range("D25").Value = '''8.58'
$ws.Range("E25").Value = '  +3.67%  '

$ws.Range("D26").Value = '''1.00'
$ws.Range("E26").Value = '  +0.34%  '

$ws.Range("E27").Value = '  +0.76%  '

$ws.Range("D28").Value = '0.0₃0778'
$ws.Range("E28").Value = '  +4.94%  '

$ws.Range("D29").Value = '''1.77'
$ws.Range("E29").Value = '  +0.96%  '

$ws.Range("D30").Value = '''170.41'
$ws.Range("E30").Value = '  -0.12%  '

$ws.Range("E31").Value = '  +0.24%  '

$ws.Range("E32").Value = '  +0.96%  '

$ws.Range("E33").Value = '  +1.18%  '

$ws.Range("E34").Value = '  -0.01%  '

$ws.Range("E35").Value = '  +3.65%  '

$ws.Range("E36").Value = '  +0.08%  '

$ws.Range("E37").Value = '  -0.62%  '

$ws.Range("E39").Value = '  +0.62%  '

$ws.Range("D40").Value = '''0.411'
$ws.Range("E40").Value = '  -0.32%  '

$ws.Range("D41").Value = '''314.55'
$ws.Range("E41").Value = '  +8.88%  '

$ws.Range("D42").Value = '''3.65'
$ws.Range("E42").Value = '  -0.85%  '

$ws.Range("D43").Value = '''138.11'
$ws.Range("E43").Value = '  -2.73%  '

$ws.Range("D45").Value = '''0.0518'
$ws.Range("E45").Value = '  -0.49%  '

$ws.Range("E46").Value = '  +2.20%  '

$ws.Range("B47").Value = 'InjectiveProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D47").Value = '''19.33'
$ws.Range("E47").Value = '  +1.72%  '

$ws.Range("B48").Value = 'Polygon'
$ws.Range("C48").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D48").Value = '''0.399'
$ws.Range("E48").Value = '  +2.48%  '

$ws.Range("E49").Value = '  +0.27%  '

$ws.Range("E50").Value = '  +0.03%  '

$ws.Range("E51").Value = '  -0.43%  '
